# Update "想去人数" (interest/attendance count) figures in column F to the
# newly scraped values, on both the "展览" sheet and the mirrored
# "全部类型" sheet (gh-pages data refresh at commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value2  = 1046
$ws1.Range("F8").Value2  = 14744
$ws1.Range("F9").Value2  = 176
$ws1.Range("F11").Value2 = 5883
$ws1.Range("F12").Value2 = 601
$ws1.Range("F13").Value2 = 83
$ws1.Range("F19").Value2 = 192
$ws1.Range("F22").Value2 = 93
$ws1.Range("F23").Value2 = 10685
$ws1.Range("F25").Value2 = 72
$ws1.Range("F26").Value2 = 107

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value2  = 1046
$ws4.Range("F9").Value2  = 14744
$ws4.Range("F10").Value2 = 176
$ws4.Range("F12").Value2 = 5883
$ws4.Range("F13").Value2 = 601
$ws4.Range("F14").Value2 = 83
$ws4.Range("F20").Value2 = 192
$ws4.Range("F23").Value2 = 93
$ws4.Range("F25").Value2 = 10685
$ws4.Range("F27").Value2 = 72
$ws4.Range("F28").Value2 = 107
